$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111356762
$ws.Range("B2").Value = 89686
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 658
$ws.Range("F2").Value = 'Rosenticka'
$ws.Range("G2").Value = 'Rhodofomes roseus'
$ws.Range("H2").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q2").Value = 553951.9614282879
$ws.Range("R2").Value = 7002044.904499132
$ws.Range("Z2").Value = '20:49'
$ws.Range("AB2").Value = '20:49'

# Row 3
$ws.Range("A3").Value = 111357157
$ws.Range("B3").Value = 78578
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = 'Lunglav'
$ws.Range("G3").Value = 'Lobaria pulmonaria'
$ws.Range("H3").Value = '(L.) Hoffm.'
$ws.Range("Q3").Value = 553906.6257793424
$ws.Range("R3").Value = 7001993.497915561
$ws.Range("Z3").Value = '21:05'
$ws.Range("AB3").Value = '21:05'

# Row 4
$ws.Range("A4").Value = 111357015
$ws.Range("B4").Value = 89845
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = 'Rynkskinn'
$ws.Range("G4").Value = 'Phlebia centrifuga'
$ws.Range("H4").Value = 'P.Karst.'
$ws.Range("Q4").Value = 553909.463631961
$ws.Range("R4").Value = 7002013.443953016
$ws.Range("Z4").Value = '20:49'
$ws.Range("AB4").Value = '20:49'

# Row 5
$ws.Range("A5").Value = 111357360
$ws.Range("B5").Value = 89405
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 1202
$ws.Range("F5").Value = 'Ullticka'
$ws.Range("G5").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H5").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q5").Value = 553854.7258749125
$ws.Range("R5").Value = 7001982.684500803
$ws.Range("Z5").Value = '21:05'
$ws.Range("AB5").Value = '21:05'

# Row 6
$ws.Range("A6").Value = 111357776
$ws.Range("B6").Value = 78578
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = 'Lunglav'
$ws.Range("G6").Value = 'Lobaria pulmonaria'
$ws.Range("H6").Value = '(L.) Hoffm.'
$ws.Range("Q6").Value = 553808.7819238321
$ws.Range("R6").Value = 7002131.15853373
$ws.Range("Z6").Value = '21:36'
$ws.Range("AB6").Value = '21:36'

# Row 7
$ws.Range("A7").Value = 111356314
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("Q7").Value = 554037.1883795768
$ws.Range("R7").Value = 7002120.944976788
$ws.Range("Z7").Value = '20:27'
$ws.Range("AB7").Value = '20:27'

# Row 8
$ws.Range("A8").Value = 111356263
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = 'Knärot'
$ws.Range("G8").Value = 'Goodyera repens'
$ws.Range("H8").Value = '(L.) R. Br.'
$ws.Range("Q8").Value = 554054.0600129352
$ws.Range("R8").Value = 7002113.991040959
$ws.Range("Z8").Value = '20:23'
$ws.Range("AB8").Value = '20:23'

# Row 9
$ws.Range("A9").Value = 111356256
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("Q9").Value = 554052.9808952439
$ws.Range("R9").Value = 7002124.374295473
$ws.Range("Z9").Value = '20:23'
$ws.Range("AB9").Value = '20:23'

# Row 10
$ws.Range("A10").Value = 111356587
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = 'VU'
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = 'Knärot'
$ws.Range("G10").Value = 'Goodyera repens'
$ws.Range("H10").Value = '(L.) R. Br.'
$ws.Range("Q10").Value = 553994.858156529
$ws.Range("R10").Value = 7002052.403435753
$ws.Range("Z10").Value = '20:39'
$ws.Range("AB10").Value = '20:39'

# Row 11
$ws.Range("A11").Value = 111356354
$ws.Range("B11").Value = 78578
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 6458
$ws.Range("F11").Value = 'Lunglav'
$ws.Range("G11").Value = 'Lobaria pulmonaria'
$ws.Range("H11").Value = '(L.) Hoffm.'
$ws.Range("Q11").Value = 554026.383447904
$ws.Range("R11").Value = 7002090.012868459
$ws.Range("Z11").Value = '20:27'
$ws.Range("AB11").Value = '20:27'

# Row 12
$ws.Range("A12").Value = 111357873
$ws.Range("B12").Value = 78578
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6458
$ws.Range("F12").Value = 'Lunglav'
$ws.Range("G12").Value = 'Lobaria pulmonaria'
$ws.Range("H12").Value = '(L.) Hoffm.'
$ws.Range("Q12").Value = 553818.3826172169
$ws.Range("R12").Value = 7002180.158265028
$ws.Range("Z12").Value = '21:36'
$ws.Range("AB12").Value = '21:36'

# Row 13
$ws.Range("A13").Value = 111358006
$ws.Range("B13").Value = 98446
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 222771
$ws.Range("F13").Value = 'Svart trolldruva'
$ws.Range("G13").Value = 'Actaea spicata'
$ws.Range("H13").Value = 'L.'
$ws.Range("Q13").Value = 553854.1622618367
$ws.Range("R13").Value = 7002179.849007829
$ws.Range("Z13").Value = '21:36'
$ws.Range("AB13").Value = '21:36'

# Row 14
$ws.Range("A14").Value = 111357080
$ws.Range("B14").Value = 89686
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 658
$ws.Range("F14").Value = 'Rosenticka'
$ws.Range("G14").Value = 'Rhodofomes roseus'
$ws.Range("H14").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q14").Value = 553906.6257793424
$ws.Range("R14").Value = 7001993.497915561
$ws.Range("Z14").Value = '20:49'
$ws.Range("AB14").Value = '20:49'

# Row 15
$ws.Range("A15").Value = 111357720
$ws.Range("B15").Value = 96348
$ws.Range("D15").Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = 'Knärot'
$ws.Range("G15").Value = 'Goodyera repens'
$ws.Range("H15").Value = '(L.) R. Br.'
$ws.Range("Q15").Value = 553822.8840132115
$ws.Range("R15").Value = 7002127.322982416
$ws.Range("Z15").Value = '21:36'
$ws.Range("AB15").Value = '21:36'

# Row 16
$ws.Range("A16").Value = 111356702
$ws.Range("B16").Value = 6202
$ws.Range("D16").Value = 'LC'
$ws.Range("E16").Value = 105336
$ws.Range("F16").Value = 'Vanlig flatbagge'
$ws.Range("G16").Value = 'Peltis ferruginea'
$ws.Range("H16").Value = '(Linnaeus, 1758)'
$ws.Range("Q16").Value = 553981.1551737323
$ws.Range("R16").Value = 7002032.27630965
$ws.Range("Z16").Value = '20:39'
$ws.Range("AB16").Value = '20:39'

# Row 17
$ws.Range("A17").Value = 111358027
$ws.Range("B17").Value = 98446
$ws.Range("D17").Value = 'LC'
$ws.Range("E17").Value = 222771
$ws.Range("F17").Value = 'Svart trolldruva'
$ws.Range("G17").Value = 'Actaea spicata'
$ws.Range("H17").Value = 'L.'
$ws.Range("Q17").Value = 553857.5193624865
$ws.Range("R17").Value = 7002168.599353628
$ws.Range("Z17").Value = '21:36'
$ws.Range("AB17").Value = '21:36'

# Row 18
$ws.Range("A18").Value = 111356632
$ws.Range("B18").Value = 78578
$ws.Range("D18").Value = 'NT'
$ws.Range("E18").Value = 6458
$ws.Range("F18").Value = 'Lunglav'
$ws.Range("G18").Value = 'Lobaria pulmonaria'
$ws.Range("H18").Value = '(L.) Hoffm.'
$ws.Range("Q18").Value = 553994.858156529
$ws.Range("R18").Value = 7002052.403435753
$ws.Range("Z18").Value = '20:39'
$ws.Range("AB18").Value = '20:39'
